$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B values are numeric-looking ("1") text in the source data; force
# text formatting so Excel doesn't coerce them to numbers on entry.
$ws.Range("B2:B18").NumberFormat = "@"

# Table of rows 2-18 (row 2 is rewritten with new content, rows 3-18 are
# newly added): row, STT(A), Đài truyền hình(B), Nội dung(C), video_link(E), ngày_giờ(F)
$data = @"
2	1	1	 Tạp chí Văn hóa – Văn nghệ		2024-01-14 16:00:00
3	2	1	 Thế giới chuyển động		2024-01-14 16:15:00
4	3	1	 An toàn giao thông		2024-01-14 16:29:00
5	4	1	 Phim Khi nắng thu về		2024-01-14 16:30:00
6	5	1	 An ninh Hậu Giang		2024-01-14 18:00:00
7	6	1	 VFC cánh đồng hội nhập		2024-01-14 18:15:00
8	7	1	 Tin tức Mekong	https://60acee235f4d5.streamlock.net:443/VODHGTV/definst/VIDEO/mp4:ttmk-140124.mp4/playlist.m3u8	2024-01-14 18:29:00
9	8	1	 Tiếp chuyển Thời sự VTV		2024-01-14 19:00:00
10	9	1	 Thời sự Hậu Giang – Thời tiết nông vụ	https://60acee235f4d5.streamlock.net:443/VODHGTV/definst/VIDEO/mp4:tshg-140124.mp4/playlist.m3u8	2024-01-14 19:40:00
11	10	1	 Hậu Giang trên đường phát triển		2024-01-14 20:15:00
12	11	1	 Phim tài liệu 20 năm thành lập tỉnh Hậu Giang (T10)		2024-01-14 20:30:00
13	12	1	 Phim Huynh đệ tương tàn (T18,19)		2024-01-14 20:50:00
14	13	1	 Thế giới chuyển động		2024-01-14 22:15:00
15	14	1	 VFC cánh đồng hội nhập		2024-01-14 22:30:00
16	15	1	 Chuyện tình tôi kể		2024-01-14 22:45:00
17	16	1	 Phim tài liệU		2024-01-14 23:15:00
18	17	1	 Trích đoạn sân khấu		2024-01-14 23:30:00
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.TrimEnd("`r")
    if ($line.Length -eq 0) { continue }
    $f = $line -split "`t"
    $rowNum = [int]$f[0]
    $stt = [int]$f[1]
    $dai = $f[2]
    $noidung = $f[3]
    $link = $f[4]
    $ngaygio = $f[5]

    $ws.Cells.Item($rowNum, 1).Value = $stt
    $ws.Cells.Item($rowNum, 2).Value = $dai
    $ws.Cells.Item($rowNum, 3).Value = $noidung
    $ws.Cells.Item($rowNum, 4).Value = ""
    $ws.Cells.Item($rowNum, 5).Value = $link
    $ws.Cells.Item($rowNum, 6).Value = $ngaygio
}
